$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "color" column header
$ws.Range("E1").Value = "color"

# Fill the new column with "black" for rows 2 through 12
$ws.Range("E2:E12").Value = "black"

# Update the selection to match the edited range
$ws.Range("E2:E12").Select()
